# custom accuracy + 데이터 1000개
# - Round the last data row (row 5) values to 2 decimal places ("custom accuracy")
# - Remove the extra data row (row 6)
# - Re-fit column widths to the now-shorter values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Round row 5 (columns B:AH) to 2 decimal places, in place, using the
#    values already present in the sheet (custom accuracy = 2 decimals).
for ($col = 2; $col -le 34; $col++) {
    $cur = $ws.Cells.Item(5, $col).Value2
    $ws.Cells.Item(5, $col).Value = $excel.WorksheetFunction.Round($cur, 2)
}

# 2) Delete the now-unneeded last row of data (row 6).
$ws.Rows.Item(6).Delete()

# 3) Shrink columns to fit the shorter, rounded values.
#    (ColumnWidth uses "characters"; the stored column width in the file is
#    ColumnWidth + 0.8333333333333334, so subtract that offset to land on
#    the desired stored widths of 7 / 8 / 21.)
$ws.Range("B1:D1").ColumnWidth = 6.166666666666667
$ws.Range("E1").ColumnWidth = 7.166666666666667
$ws.Range("F1:S1").ColumnWidth = 6.166666666666667
$ws.Range("T1:U1").ColumnWidth = 7.166666666666667
$ws.Range("V1:AH1").ColumnWidth = 6.166666666666667
